# Updates the cryptos list (Coin / Link / Price / Volume(1h)) for rows 2-51
# to reflect the latest scraped values, per commit:
# "Updated cryptos list on Mon Aug 14 14:26:14 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row|Coin|Link|Price|Volume(1h)
$rowData = @'
2|Bitcoin|https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc|29.379.82|  +0.01%  
3|Ethereum|https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth|1.845.84|  -0.20%  
4|TetherUSD|https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt|0.9976|  -0.16%  
5|BNB|https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb|240.18|  -0.15%  
6|XRP|https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp|0.6288|  -0.01%  
7|USDC|https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc|0.9993|  -0.11%  
8|Dogecoin|https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge|0.07458|  -2.04%  
9|Cardano|https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada|0.2906|  -0.38%  
10|Solana|https://coinranking.com/coin/zNZHO_Sjf+solana-sol|24.49|  -0.50%  
11|TRON|https://coinranking.com/coin/qUhEFk1I61atv+tron-trx|0.07739|  -0.16%  
12|WrappedEther|https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth|1.841.98|  -2.45%  
13|Polkadot|https://coinranking.com/coin/25W7FG7om+polkadot-dot|4.995|  -0.62%  
14|Polygon|https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic|0.6782|  -0.38%  
15|ShibaInu|https://coinranking.com/coin/xz24e0BjL+shibainu-shib|0.00001044|  -0.82%  
16|Litecoin|https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc|82.12|  -1.29%  
17|Uniswap|https://coinranking.com/coin/_H5FVG9iW+uniswap-uni|6.194|  +1.00%  
18|WrappedBTC|https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc|29.408.31|  +0.08%  
19|BitcoinCash|https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch|228.20|  -0.54%  
20|Avalanche|https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax|12.32|  -0.12%  
21|Dai|https://coinranking.com/coin/MoTuySvg7+dai-dai|0.9996|  -0.05%  
22|Chainlink|https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link|7.506|  +0.44%  
23|BinanceUSD|https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd|0.9999|  -0.11%  
24|Monero|https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr|159.19|  +0.33%  
25|Cosmos|https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom|8.500|  +0.60%  
26|Stellar|https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm|0.1367|  -1.68%  
27|EthereumClassic|https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc|17.51|  -1.10%  
28|Hedera|https://coinranking.com/coin/jad286TjB+hedera-hbar|0.06451|  +14.86%  
29|Toncoin|https://coinranking.com/coin/67YlI0K1b+toncoin-ton|1.417|  -2.38%  
30|PancakeSwap|https://coinranking.com/coin/ncYFcP709+pancakeswap-cake|1.484|  +0.72%  
31|Filecoin|https://coinranking.com/coin/ymQub4fuB+filecoin-fil|4.088|  -0.62%  
32|InternetComputer(DFINITY)|https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp|4.092|  +0.70%  
33|LidoDAOToken|https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo|1.834|  +0.05%  
34|ARBITRUM|https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb|1.142|  -1.44%  
35|ImmutableX|https://coinranking.com/coin/Z96jIvLU7+immutablex-imx|0.6951|  -0.57%  
36|HuobiToken|https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht|2.584|  -0.11%  
37|Maker|https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr|1.261.20|  +1.97%  
38|VeChain|https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet|0.01845|  +2.24%  
39|MXToken|https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx|2.839|  +4.00%  
40|FraxShare|https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs|6.790|  +5.61%  
41|TrustWalletToken|https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt|0.9328|  +3.33%  
42|PaxDollar|https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp|0.9990|  -0.11%  
43|RocketPoolETH|https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth|2.007.00|  +1.28%  
44|Quant|https://coinranking.com/coin/bauj_21eYVwso+quant-qnt|101.31|  -0.33%  
45|Aave|https://coinranking.com/coin/ixgUfzmLR+aave-aave|65.98|  +0.49%  
46|RenderToken|https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr|1.734|  +2.94%  
47|Aptos|https://coinranking.com/coin/HGYj5JCv5+aptos-apt|7.080|  -1.41%  
48|Algorand|https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo|0.1157|  +0.21%  
49|EnergySwap|https://coinranking.com/coin/SbWqqTui-+energyswap-ens|9.027|  +0.40%  
50|TheSandbox|https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand|0.3941|  -1.43%  
51|BabyDogeCoin|https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge|0.00000000114|  +2.12%  
'@

$lines = $rowData -split "`n"
foreach ($line in $lines) {
    $line = $line.TrimEnd("`r")
    if ($line.Length -eq 0) { continue }
    $parts = $line -split '\|', 5
    $r = [int]$parts[0]
    $coin = $parts[1]
    $link = $parts[2]
    $price = $parts[3]
    $volume = $parts[4]

    $ws.Cells.Item($r, 2).Value = $coin
    $ws.Cells.Item($r, 3).Value = $link

    # Price column holds plain text (e.g. "0.9976", "29.379.82"); force
    # text format first so Excel does not coerce it into a Double and
    # strip formatting (trailing zeros, multi-dot groupings, etc.), then
    # restore the default "Normal" style so no stray number format is
    # left behind on the cell.
    $priceCell = $ws.Cells.Item($r, 4)
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $price
    $priceCell.Style = "Normal"

    $ws.Cells.Item($r, 5).Value = $volume
}
